$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 9260042
$ws.Range("I43").Value = 833.3333
$ws.Range("J43").Value = 18519252
$ws.Range("K43").Value = 833.3333
$ws.Range("L43").Value = 18519252
$ws.Range("M43").Value = -764.3333
$ws.Range("N43").Value = -18519390

$ws.Range("H53").Value = 1393.2
$ws.Range("I53").Value = 1926.7142
$ws.Range("K53").Value = 1926.7142
$ws.Range("M53").Value = -1289.7142

$ws.Range("H76").Value = 4812.5
$ws.Range("I76").Value = 6500
$ws.Range("K76").Value = 6500
$ws.Range("M76").Value = -6185

$ws.Range("H79").Value = 4812.5
$ws.Range("I79").Value = 6500
$ws.Range("K79").Value = 6500
$ws.Range("M79").Value = -5408

$ws.Range("H125").Value = 805
$ws.Range("I125").Value = 760
$ws.Range("J125").Value = 850
$ws.Range("K125").Value = 6840
$ws.Range("L125").Value = 7650
$ws.Range("M125").Value = -4380
$ws.Range("N125").Value = -12570

$ws.Range("H132").Value = 7534.0527
$ws.Range("I132").Value = 4395.5835
$ws.Range("J132").Value = 12914.286
$ws.Range("K132").Value = 13186.7505
$ws.Range("L132").Value = 38742.858
$ws.Range("M132").Value = -10656.7505
$ws.Range("N132").Value = -43802.858

$ws.Range("H133").Value = 34499.5
$ws.Range("J133").Value = 34499.5
$ws.Range("L133").Value = 34499.5
$ws.Range("N133").Value = -44619.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5770.1787
$ws.Range("I32").Value = 5692.731
$ws.Range("J32").Value = 6777
$ws.Range("K32").Value = 5692.731
$ws.Range("L32").Value = 6777
$ws.Range("M32").Value = -5405.731
$ws.Range("N32").Value = -7351

$ws.Range("H61").Value = 1209.5714
$ws.Range("I61").Value = 976.8889
$ws.Range("J61").Value = 1628.4
$ws.Range("K61").Value = 976.8889
$ws.Range("L61").Value = 1628.4
$ws.Range("M61").Value = -764.8889
$ws.Range("N61").Value = -2052.4

$ws.Range("H122").Value = 1417.6
$ws.Range("I122").Value = 1544
$ws.Range("J122").Value = 1333.3334
$ws.Range("K122").Value = 4632
$ws.Range("L122").Value = 4000.0002
$ws.Range("M122").Value = -2182
$ws.Range("N122").Value = -8900.0002

$ws.Range("H123").Value = 75000
$ws.Range("J123").Value = 75000
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -84800

$ws.Range("H132").Value = 2597.4211
$ws.Range("I132").Value = 2551.8845
$ws.Range("J132").Value = 2696.0833
$ws.Range("K132").Value = 7655.6535
$ws.Range("L132").Value = 8088.249899999999
$ws.Range("M132").Value = -5125.6535
$ws.Range("N132").Value = -13148.2499

$ws.Range("H136").Value = 1209.5714
$ws.Range("I136").Value = 976.8889
$ws.Range("J136").Value = 1628.4
$ws.Range("K136").Value = 2930.6667
$ws.Range("L136").Value = 4885.200000000001
$ws.Range("M136").Value = -380.6667000000002
$ws.Range("N136").Value = -9985.200000000001

$ws.Range("H137").Value = 85500
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 85500
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 85500
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -95700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 71430230
$ws.Range("I99").Value = 83334950
$ws.Range("J99").Value = 1900
$ws.Range("K99").Value = 83334950
$ws.Range("L99").Value = 1900
$ws.Range("M99").Value = -83333452
$ws.Range("N99").Value = -4896

$ws.Range("H132").Value = 2021467.4
$ws.Range("J132").Value = 2021467.4
$ws.Range("L132").Value = 2021467.4
$ws.Range("N132").Value = -2031587.4

$ws.Range("H134").Value = 5570.2607
$ws.Range("I134").Value = 1196.0476
$ws.Range("K134").Value = 3588.142800000001
$ws.Range("M134").Value = -1053.142800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 18000
$ws.Range("J48").Value = 18000
$ws.Range("L48").Value = 18000
$ws.Range("N48").Value = -18952

$ws.Range("H62").Value = 4764409.5
$ws.Range("I62").Value = 2565.9756
$ws.Range("K62").Value = 2565.9756
$ws.Range("M62").Value = -1941.9756

$ws.Range("H65").Value = 4764409.5
$ws.Range("I65").Value = 2565.9756
$ws.Range("K65").Value = 12829.878
$ws.Range("M65").Value = -9709.878000000001

$ws.Range("H99").Value = 1431.7059
$ws.Range("I99").Value = 1632.2858
$ws.Range("J99").Value = 1291.3
$ws.Range("K99").Value = 1632.2858
$ws.Range("L99").Value = 1291.3
$ws.Range("M99").Value = -134.2858000000001
$ws.Range("N99").Value = -4287.3

$ws.Range("H126").Value = 1431.7059
$ws.Range("I126").Value = 1632.2858
$ws.Range("J126").Value = 1291.3
$ws.Range("K126").Value = 4896.857400000001
$ws.Range("L126").Value = 3873.9
$ws.Range("M126").Value = -2426.857400000001
$ws.Range("N126").Value = -8813.9

$ws.Range("H132").Value = 4577.7646
$ws.Range("I132").Value = 4878.5386
$ws.Range("J132").Value = 3600.25
$ws.Range("K132").Value = 14635.6158
$ws.Range("L132").Value = 10800.75
$ws.Range("M132").Value = -12105.6158
$ws.Range("N132").Value = -15860.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1601.0769
$ws.Range("I136").Value = 903
$ws.Range("J136").Value = 2199.4285
$ws.Range("K136").Value = 2709
$ws.Range("L136").Value = 6598.2855
$ws.Range("M136").Value = 2391
$ws.Range("N136").Value = -16798.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3106.4211
$ws.Range("I122").Value = 2002.75
$ws.Range("J122").Value = 3909.0908
$ws.Range("K122").Value = 6008.25
$ws.Range("L122").Value = 11727.2724
$ws.Range("M122").Value = -3558.25
$ws.Range("N122").Value = -16627.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1487.1333
$ws.Range("I7").Value = 1343
$ws.Range("K7").Value = 1343
$ws.Range("M7").Value = -1231

$ws.Range("H16").Value = 1156.1666
$ws.Range("I16").Value = 887.7273
$ws.Range("K16").Value = 887.7273
$ws.Range("M16").Value = -717.7273

$ws.Range("H40").Value = 2252.6924
$ws.Range("I40").Value = 2023.75
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 2023.75
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -1887.75
$ws.Range("N40").Value = -5272

$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21498

$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -67488

$ws.Range("H100").Value = 2572.5715
$ws.Range("I100").Value = 2001.5
$ws.Range("K100").Value = 2001.5
$ws.Range("M100").Value = -1460.5

$ws.Range("H122").Value = 13892276
$ws.Range("I122").Value = 35716784
$ws.Range("J122").Value = 3951.0908
$ws.Range("K122").Value = 107150352
$ws.Range("L122").Value = 11853.2724
$ws.Range("M122").Value = -107147902
$ws.Range("N122").Value = -16753.2724

$ws.Range("H126").Value = 1487.1333
$ws.Range("I126").Value = 1343
$ws.Range("K126").Value = 4029
$ws.Range("M126").Value = -1559

$ws.Range("H132").Value = 61294.47
$ws.Range("I132").Value = 2571.1428
$ws.Range("J132").Value = 102400.8
$ws.Range("K132").Value = 7713.428400000001
$ws.Range("L132").Value = 307202.4
$ws.Range("M132").Value = -5183.428400000001
$ws.Range("N132").Value = -312262.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2165.0833
$ws.Range("I81").Value = 1548.1
$ws.Range("K81").Value = 3096.2
$ws.Range("M81").Value = -2035.2

$ws.Range("H84").Value = 2165.0833
$ws.Range("I84").Value = 1548.1
$ws.Range("K84").Value = 15481
$ws.Range("M84").Value = -10177

$ws.Range("H94").Value = 29500
$ws.Range("J94").Value = 29500
$ws.Range("L94").Value = 29500
$ws.Range("N94").Value = -31302

$ws.Range("H100").Value = 944.8889
$ws.Range("I100").Value = 850.6667
$ws.Range("J100").Value = 1133.3334
$ws.Range("K100").Value = 1701.3334
$ws.Range("L100").Value = 2266.6668
$ws.Range("M100").Value = -1160.3334
$ws.Range("N100").Value = -3348.6668

$ws.Range("H128").Value = 99990
$ws.Range("J128").Value = 99990
$ws.Range("L128").Value = 99990
$ws.Range("N128").Value = -109950

$ws.Range("H139").Value = 34714
$ws.Range("I139").Value = 34714
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 34714
$ws.Range("L139").ClearContents()
$ws.Range("M139").Value = -29574
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 34698
$ws.Range("J140").Value = 34698
$ws.Range("L140").Value = 34698
$ws.Range("N140").Value = -45058

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
